$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 2.051623430250302
$ws.Range("C3").Value2 = 0.5608039820102064
$ws.Range("C4").Value2 = 0.3111664295874803
$ws.Range("C5").Value2 = 0.04768495693301134
$ws.Range("C6").Value2 = -0.001699323063178326
$ws.Range("C7").Value2 = -0.1434825472076895
$ws.Range("C8").Value2 = -0.2622913429601894
$ws.Range("C9").Value2 = -0.5027825232064487
$ws.Range("C10").Value2 = -0.6086099171784909
